$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: plain "NNN.NNN" decimal strings would be auto-converted to numbers by Excel,
# losing the exact text formatting (trailing zeros, float precision). Such values are entered
# with a leading apostrophe to force them to stay text, matching the original text-cell content.
$numericLike = '^\d+\.\d+$'

# Row 2
$dVal = "60.149.51"
if ($dVal -match $numericLike) { $ws.Range("D2").Value = "'" + $dVal } else { $ws.Range("D2").Value = $dVal }
$ws.Range("E2").Value = "  +1.87%  "

# Row 3
$dVal = "2.668.48"
if ($dVal -match $numericLike) { $ws.Range("D3").Value = "'" + $dVal } else { $ws.Range("D3").Value = $dVal }
$ws.Range("E3").Value = "  +0.32%  "

# Row 4
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$dVal = "520.00"
if ($dVal -match $numericLike) { $ws.Range("D5").Value = "'" + $dVal } else { $ws.Range("D5").Value = $dVal }
$ws.Range("E5").Value = "  +1.12%  "

# Row 6
$dVal = "145.43"
if ($dVal -match $numericLike) { $ws.Range("D6").Value = "'" + $dVal } else { $ws.Range("D6").Value = $dVal }
$ws.Range("E6").Value = "  +0.83%  "

# Row 7
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$dVal = "0.578"
if ($dVal -match $numericLike) { $ws.Range("D8").Value = "'" + $dVal } else { $ws.Range("D8").Value = $dVal }
$ws.Range("E8").Value = "  +1.43%  "

# Row 9
$dVal = "2.676.22"
if ($dVal -match $numericLike) { $ws.Range("D9").Value = "'" + $dVal } else { $ws.Range("D9").Value = $dVal }
$ws.Range("E9").Value = "  +0.68%  "

# Row 10
$dVal = "6.43"
if ($dVal -match $numericLike) { $ws.Range("D10").Value = "'" + $dVal } else { $ws.Range("D10").Value = $dVal }
$ws.Range("E10").Value = "  +1.68%  "

# Row 11
$ws.Range("E11").Value = "  -0.63%  "

# Row 12
$dVal = "0.339"
if ($dVal -match $numericLike) { $ws.Range("D12").Value = "'" + $dVal } else { $ws.Range("D12").Value = $dVal }
$ws.Range("E12").Value = "  +0.74%  "

# Row 13
$ws.Range("E13").Value = "  +1.48%  "

# Row 14
$dVal = "3.138.86"
if ($dVal -match $numericLike) { $ws.Range("D14").Value = "'" + $dVal } else { $ws.Range("D14").Value = $dVal }
$ws.Range("E14").Value = "  +1.31%  "

# Row 15
$dVal = "60.151.85"
if ($dVal -match $numericLike) { $ws.Range("D15").Value = "'" + $dVal } else { $ws.Range("D15").Value = $dVal }
$ws.Range("E15").Value = "  +1.93%  "

# Row 16
$dVal = "21.22"
if ($dVal -match $numericLike) { $ws.Range("D16").Value = "'" + $dVal } else { $ws.Range("D16").Value = $dVal }
$ws.Range("E16").Value = "  +0.66%  "

# Row 17
$dVal = "2.771.88"
if ($dVal -match $numericLike) { $ws.Range("D17").Value = "'" + $dVal } else { $ws.Range("D17").Value = $dVal }
$ws.Range("E17").Value = "  +4.57%  "

# Row 18
$ws.Range("E18").Value = "  +0.76%  "

# Row 19
$dVal = "349.03"
if ($dVal -match $numericLike) { $ws.Range("D19").Value = "'" + $dVal } else { $ws.Range("D19").Value = $dVal }
$ws.Range("E19").Value = "  +1.45%  "

# Row 20
$dVal = "4.54"
if ($dVal -match $numericLike) { $ws.Range("D20").Value = "'" + $dVal } else { $ws.Range("D20").Value = $dVal }
$ws.Range("E20").Value = "  -0.15%  "

# Row 21
$dVal = "10.50"
if ($dVal -match $numericLike) { $ws.Range("D21").Value = "'" + $dVal } else { $ws.Range("D21").Value = $dVal }
$ws.Range("E21").Value = "  +1.32%  "

# Row 22
$dVal = "6.29"
if ($dVal -match $numericLike) { $ws.Range("D22").Value = "'" + $dVal } else { $ws.Range("D22").Value = $dVal }
$ws.Range("E22").Value = "  +3.16%  "

# Row 23
$dVal = "0.999"
if ($dVal -match $numericLike) { $ws.Range("D23").Value = "'" + $dVal } else { $ws.Range("D23").Value = $dVal }
$ws.Range("E23").Value = "  +0.24%  "

# Row 24
$dVal = "62.66"
if ($dVal -match $numericLike) { $ws.Range("D24").Value = "'" + $dVal } else { $ws.Range("D24").Value = $dVal }
$ws.Range("E24").Value = "  +2.66%  "

# Row 25
$dVal = "0.420"
if ($dVal -match $numericLike) { $ws.Range("D25").Value = "'" + $dVal } else { $ws.Range("D25").Value = $dVal }
$ws.Range("E25").Value = "  -0.14%  "

# Row 26
$dVal = "0.167"
if ($dVal -match $numericLike) { $ws.Range("D26").Value = "'" + $dVal } else { $ws.Range("D26").Value = $dVal }
$ws.Range("E26").Value = "  +4.31%  "

# Row 27
$dVal = "0.995"
if ($dVal -match $numericLike) { $ws.Range("D27").Value = "'" + $dVal } else { $ws.Range("D27").Value = $dVal }
$ws.Range("E27").Value = "  +0.45%  "

# Row 28
$dVal = "0.0₃0809"
if ($dVal -match $numericLike) { $ws.Range("D28").Value = "'" + $dVal } else { $ws.Range("D28").Value = $dVal }
$ws.Range("E28").Value = "  +0.10%  "

# Row 29
$dVal = "7.22"
if ($dVal -match $numericLike) { $ws.Range("D29").Value = "'" + $dVal } else { $ws.Range("D29").Value = $dVal }
$ws.Range("E29").Value = "  +1.44%  "

# Row 30
$dVal = "6.83"
if ($dVal -match $numericLike) { $ws.Range("D30").Value = "'" + $dVal } else { $ws.Range("D30").Value = $dVal }
$ws.Range("E30").Value = "  +6.20%  "

# Row 31
$ws.Range("E31").Value = "  +0.18%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$dVal = "1.59"
if ($dVal -match $numericLike) { $ws.Range("D32").Value = "'" + $dVal } else { $ws.Range("D32").Value = $dVal }
$ws.Range("E32").Value = "  +0.64%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$dVal = "19.00"
if ($dVal -match $numericLike) { $ws.Range("D33").Value = "'" + $dVal } else { $ws.Range("D33").Value = $dVal }
$ws.Range("E33").Value = "  +0.54%  "

# Row 34
$dVal = "148.34"
if ($dVal -match $numericLike) { $ws.Range("D34").Value = "'" + $dVal } else { $ws.Range("D34").Value = $dVal }
$ws.Range("E34").Value = "  -1.09%  "

# Row 35
$dVal = "4.29"
if ($dVal -match $numericLike) { $ws.Range("D35").Value = "'" + $dVal } else { $ws.Range("D35").Value = $dVal }
$ws.Range("E35").Value = "  +6.36%  "

# Row 36
$dVal = "0.947"
if ($dVal -match $numericLike) { $ws.Range("D36").Value = "'" + $dVal } else { $ws.Range("D36").Value = $dVal }
$ws.Range("E36").Value = "  -6.19%  "

# Row 37
$ws.Range("E37").Value = "  +5.78%  "

# Row 38
$ws.Range("E38").Value = "  +10.41%  "

# Row 39
$dVal = "0.868"
if ($dVal -match $numericLike) { $ws.Range("D39").Value = "'" + $dVal } else { $ws.Range("D39").Value = $dVal }
$ws.Range("E39").Value = "  +1.48%  "

# Row 40
$dVal = "36.67"
if ($dVal -match $numericLike) { $ws.Range("D40").Value = "'" + $dVal } else { $ws.Range("D40").Value = $dVal }

# Row 41
$dVal = "3.68"
if ($dVal -match $numericLike) { $ws.Range("D41").Value = "'" + $dVal } else { $ws.Range("D41").Value = $dVal }
$ws.Range("E41").Value = "  -0.39%  "

# Row 42
$dVal = "280.62"
if ($dVal -match $numericLike) { $ws.Range("D42").Value = "'" + $dVal } else { $ws.Range("D42").Value = $dVal }
$ws.Range("E42").Value = "  -0.26%  "

# Row 43
$dVal = "0.0988"
if ($dVal -match $numericLike) { $ws.Range("D43").Value = "'" + $dVal } else { $ws.Range("D43").Value = $dVal }
$ws.Range("E43").Value = "  +0.31%  "

# Row 44
$dVal = "0.995"
if ($dVal -match $numericLike) { $ws.Range("D44").Value = "'" + $dVal } else { $ws.Range("D44").Value = $dVal }
$ws.Range("E44").Value = "  -0.23%  "

# Row 45
$dVal = "19.92"
if ($dVal -match $numericLike) { $ws.Range("D45").Value = "'" + $dVal } else { $ws.Range("D45").Value = $dVal }
$ws.Range("E45").Value = "  +1.93%  "

# Row 46
$dVal = "0.606"
if ($dVal -match $numericLike) { $ws.Range("D46").Value = "'" + $dVal } else { $ws.Range("D46").Value = $dVal }
$ws.Range("E46").Value = "  -1.20%  "

# Row 47
$dVal = "2.119.11"
if ($dVal -match $numericLike) { $ws.Range("D47").Value = "'" + $dVal } else { $ws.Range("D47").Value = $dVal }
$ws.Range("E47").Value = "  +6.87%  "

# Row 48
$dVal = "0.0539"
if ($dVal -match $numericLike) { $ws.Range("D48").Value = "'" + $dVal } else { $ws.Range("D48").Value = $dVal }
$ws.Range("E48").Value = "  +0.54%  "

# Row 49
$dVal = "4.83"
if ($dVal -match $numericLike) { $ws.Range("D49").Value = "'" + $dVal } else { $ws.Range("D49").Value = $dVal }
$ws.Range("E49").Value = "  +2.54%  "

# Row 50
$ws.Range("E50").Value = "  +2.10%  "

# Row 51
$dVal = "10.44"
if ($dVal -match $numericLike) { $ws.Range("D51").Value = "'" + $dVal } else { $ws.Range("D51").Value = $dVal }
$ws.Range("E51").Value = "  +1.49%  "
